# Generate Report for Handback
#
# - Updates the "Ready for handoff" status text everywhere to
#   "Handed back: in sync with en-US" (Overview + zh-cn + de-de sheets,
#   which all share the same string).
# - Fills in the "Latest Target File" / "Latest Handback File" /
#   "Latest Handback DateTime" columns (I/J/K) on the zh-cn and de-de
#   sheets for both data rows, turning column I into a hyperlink to the
#   corresponding source .md file (mirroring column A's hyperlink).
# - Widens columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d9ec3d7d10633a279f3767feec83bedf118181fe/e2e/"

$file1 = "6736bf82-0b4d-4c9a-8e94-b430845ec0ba.md"
$file2 = "ecf389f6-1002-4fb6-8099-c2d03c5786c2.md"

function Update-StatusOnSheet($ws, [string[]]$cols) {
    foreach ($col in $cols) {
        foreach ($row in 2..3) {
            $rng = $ws.Range($col + $row)
            if ($rng.Text -eq $oldStatus) {
                $rng.Value = $newStatus
            }
        }
    }
}

# --- Overview sheet: E/F hold the per-language status ---
$wsOverview = $wb.Worksheets.Item("Overview")
Update-StatusOnSheet $wsOverview @("E", "F")
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# --- zh-cn / de-de language sheets ---
$langSheets = @(
    @{
        Name = "zh-cn"
        Xlf1 = "6736bf82-0b4d-4c9a-8e94-b430845ec0ba.64c60452133987ac606d958a2518bd3fa0f6155b.zh-cn.xlf"
        Xlf2 = "ecf389f6-1002-4fb6-8099-c2d03c5786c2.700323a7a3499e645250cd61dcacde46af7349f7.zh-cn.xlf"
        DateTime = "2016-09-01 14:37:55"
    },
    @{
        Name = "de-de"
        Xlf1 = "6736bf82-0b4d-4c9a-8e94-b430845ec0ba.64c60452133987ac606d958a2518bd3fa0f6155b.de-de.xlf"
        Xlf2 = "ecf389f6-1002-4fb6-8099-c2d03c5786c2.700323a7a3499e645250cd61dcacde46af7349f7.de-de.xlf"
        DateTime = "2016-09-01 14:38:14"
    }
)

foreach ($info in $langSheets) {
    $ws = $wb.Worksheets.Item($info.Name)

    Update-StatusOnSheet $ws @("C")

    # Row 2 -> file1, Row 3 -> file2
    $ws.Range("J2").Value = $info.Xlf1
    $ws.Range("K2").Value = $info.DateTime
    $ws.Range("J3").Value = $info.Xlf2
    $ws.Range("K3").Value = $info.DateTime

    $ws.Hyperlinks.Add($ws.Range("I2"), ($baseUrl + $file1), "", "", $file1)
    $ws.Hyperlinks.Add($ws.Range("I3"), ($baseUrl + $file2), "", "", $file2)

    $ws.Columns.Item(3).ColumnWidth = 29.166666666666668
    $ws.Columns.Item(9).ColumnWidth = 39.166666666666664
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664
}

Write-Output "Handback report generated."
